$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - new requisito: Html e Css da pagina da organizacao
$ws.Range("B12").Value = "Html e Css da página da organização"
$ws.Range("C12").Value = "Realizar toda a página da organização em html e css"
$ws.Range("D12").Value = 43785

# Row 13 - new requisito: Criacao da funcionalidade Filtragem
$ws.Range("B13").Value = "Criação da funcionalidade Filtragem"
$ws.Range("C13").Value = "Poder filtrar entre pedidos pagos, prontos, entregues"
$ws.Range("D13").Value = 43801

# Apply the date format / centered alignment to the first new date cell,
# then copy that exact format onto the second so both share one style entry
# (instead of each creating its own duplicate style).
$ws.Range("D12").NumberFormat = "mm-dd-yy"
$ws.Range("D12").HorizontalAlignment = -4108
$ws.Range("D12").Copy()
$ws.Range("D13").PasteSpecial(-4122)

# Match formatting (border/alignment) of column B on these new rows for column C
$ws.Range("B12").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# Row 14 - C14 gains an (empty) cell entry alongside B14, matching B14's formatting
$ws.Range("B14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# Move the active selection down to B14
$ws.Range("B14").Select()
